$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A8").Value = "Controlador"
$ws.Range("B8").Value = "Controller"

$ws.Range("A9").Select()
